$wb = $excel.ActiveWorkbook

# The new weekly ranking sheet starts life as a copy of the most recent
# week's sheet (same header/styles), inserted right after it, then the
# rank/title/author/latest_episode values are overwritten with this
# week's data.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2025-10-20"

$rows = @(
    @(1, "ワンパンマン", "原作/ＯＮＥ 作画/村田雄介", "213撃目"),
    @(2, "転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～", "zunta(作画) はらわたさいぞう(原作)", "第33話：誰にも負けない完璧②"),
    @(3, "悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる", "こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)", "第18話-1：「R18」"),
    @(4, "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―", "光永康則", "第７０話「突貫停止」②"),
    @(5, "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。", "マツモトケンゴ", "第６６話　登山の戦いが始まった（１）"),
    @(6, "帰ってください！ 阿久津さん", "長岡太一(著者)", "番外編㉓"),
    @(7, "王子様の友達", "すけろく(著者)", "番外編【マンガ総選挙１位・マニフェスト実施】"),
    @(8, "元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～", "沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)", "第79話その2"),
    @(9, "実は俺、最強でした？", "原作：澄守 彩 漫画：高橋 愛", "第129話　未知なる脅威!?"),
    @(10, "辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？", "tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)", "第11話-2"),
    @(11, "金属スライムを倒しまくった俺が【黒鋼の王】と呼ばれるまで", "藤屋いずこ(著者) 温泉カピバラ(原作) 山椒魚(キャラクター原案)", "第15章-2"),
    @(12, "怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった", "菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)", "第14話"),
    @(13, "異世界魔王と召喚少女の奴隷魔術", "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大", "第129話　三人の門出（後編）"),
    @(14, "クセ強彼女は床にいざなう", "須河篤志(著者)", "第16話前半"),
    @(15, "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～", "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり", "第５３話　虎の尾を踏む器用貧乏（４）"),
    @(16, "ゲーム世界で魔物に転生してしまった俺、前世で推しだったヒロインを拾ってしまう", "三部べべ(漫画) ねうしとら(原作)", "第3話-1"),
    @(17, "ダークサモナーとデキている", "車王(著者)", "第78話"),
    @(18, "異世界のんびり農家", "剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)", "第308話"),
    @(19, "まんきつしたい常連さん", "しんみりん(著者)", "第49話前編"),
    @(20, "リビルドワールド", "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)", "第74話➁"),
    @(21, "【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！", "島知宏 音速炒飯 有都あらゆる", "討伐のお礼ですわ！"),
    @(22, "望まぬ不死の冒険者", "中曽根ハイジ（漫画） 丘野 優（原作） じゃいあん（キャラクター原案）", "第61話　早く捨てろ"),
    @(23, "剥かせて！竜ケ崎さん", "一智和智", "大学生編 第15話"),
    @(24, "転生貴族の異世界冒険録 ～自重を知らない神々の使徒～", "夜州 nini 藻", "第70話(前編)"),
    @(25, "独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～", "漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき", "第35話 独身貴族はバーでハイボールを作る（1）"),
    @(26, "世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜", "戸賀 環 坂木持丸 riritto", "第54話②　闇の女神と仲良くなってみた"),
    @(27, "よくわからないけれど異世界に転生していたようです", "内々けやき あし カオミン", "第141話 よくわからないけれど超理論が生み出されたようです（１）"),
    @(28, "聖者無双", "漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime", "第93話　妥協ライン（後半）"),
    @(29, "ライドンキング", "馬場康誌", "第84話 大統領と賢者の過ち（前編）"),
    @(30, "ぽんドロイド！ はまさん", "はれやまはれぞう(著者)", "第10話"),
    @(31, "魔導具師ダリヤはうつむかない ～Dahliya Wilts No More～", "漫画：住川惠 原作：甘岸久弥(｢魔導具師ダリヤはうつむかない ～今日から自由な職人ライフ～｣MFブックス刊) キャラクター原案：景、駒田ハチ", "第48話 緑の塔夏祭り夕食会②"),
    @(32, "小林さんちのメイドラゴン", "クール教信者", "第152話"),
    @(33, "願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜", "ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)", "第7話-2：新しい目標"),
    @(34, "追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。", "六志麻あさ 業務用餅 kisui", "第７４話ー①"),
    @(35, "勇者に全部奪われた俺は勇者の母親とパーティを組みました！", "久遠まこと(著者) 石のやっさん(原作)", "コミックス６巻発売告知"),
    @(36, "理想のヒモ生活", "日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)", "第88話　その3"),
    @(37, "姫様“拷問”の時間です", "原作:春原ロビンソン　漫画:ひらけい", "拷問153"),
    @(38, "地味子な三葉さんが僕を誘惑する", "はぶらえる(著者)", "第11話後半"),
    @(39, "宇崎ちゃんは遊びたい！", "丈(著者)", "第128話"),
    @(40, "俺は星間国家の悪徳領主！", "灘島かい（漫画） 三嶋与夢（原作） 高峰ナダレ（キャラクター原案）", "第41話　究極にして至高（前編）"),
    @(41, "ライブダンジョン！", "ことりりょう(作画) dy冷凍(原作) Mika Pikazo(キャラクター原案)", "第90話前半"),
    @(42, "落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～", "村上よしゆき 茨木野 あるてら", "【描き下ろしイラスト】休載です！（その２）"),
    @(43, "ギャルとダンジョンと周回遅れの探索英雄譚", "漫画家： 水田ケンジ 原作：榊一郎 キャラクター原案：黒獅子", "第4話"),
    @(44, "異世界でも無難に生きたい症候群", "原作：安泰（一二三書房刊） 漫画：笹峰コウ キャラクター原案：ひたきゆう", "第32話①"),
    @(45, "アザミヤコを好きになる", "ユニティコング(原作) ツノニガウ(作画)", "コミックス第2巻情報公開!!!"),
    @(46, "俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。", "狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)", "第26話-2"),
    @(47, "魔石グルメ　魔物の力を食べたオレは最強！", "菅原健二(作画) 結城涼(原作) 成瀬ちさと(キャラクター原案)", "第68話後半"),
    @(48, "ゲーム悪役貴族に転生した俺は、チート筋肉で無双する", "昼行燈（原作） しいたけ元帥（漫画）", "第31話"),
    @(49, "めっちゃ召喚された件 THE COMIC", "漫画：六甲島カモメ 原作：さいとうさ キャラクター原案：ツグトク", "第49話①"),
    @(50, "最弱貴族に転生したので悪役たちを集めてみた", "空野進 sorani ファルまろ", "第13話　最弱貴族、悪役令嬢に賭けを挑む（２）")
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
